$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 12 ("8. bs200 hu300") results
$ws.Range("B12").Value = 54.152000000000001
$ws.Range("D12").Value = 54.152000000000001
$ws.Range("E12").Value = 20
$ws.Range("F12").Value = 200
$ws.Range("G12").Value = 300

# Add MAX formulas to the totals row
$ws.Range("B22").Formula = "=MAX(B5:B21)"
$ws.Range("D22").Formula = "=MAX(D5:D21)"

# Update the active selection to B13, matching the saved view state
$ws.Range("B13").Select()
